$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.105.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.107.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("E9").Value = '  +1.53%  '
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.384'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.641.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("E14").Value = '  +4.74%  '
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.110.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.102.90'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("E18").Value = '  +2.18%  '
$ws.Range("E19").Value = '  +0.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.40%  '
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("E23").Value = '  +2.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0934'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.22%  '
$ws.Range("E30").Value = '  +2.74%  '
$ws.Range("E31").Value = '  +4.01%  '
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '155.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("E34").Value = '  +2.95%  '
$ws.Range("E35").Value = '  +5.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.76%  '
$ws.Range("E37").Value = '  +5.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0689'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("E39").Value = '  +3.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.150.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.82'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  +6.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.301.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = '  +2.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.970'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.760'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '263.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.93%  '
